$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -5
$ws.Range("F5").Value = -1
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = 0
$ws.Range("F16").Value = -6
$ws.Range("F22").Value = -8
$ws.Range("F23").Value = -5
$ws.Range("F24").Value = 2
$ws.Range("F25").Value = -5
$ws.Range("F34").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("F50").Value = 2
$ws.Range("F51").Value = -1
$ws.Range("F55").Value = 6
